$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 549.5  # was 466.33334
$ws.Range("I18").Value = 549.5  # was 466.33334
$ws.Range("K18").Value = 549.5  # was 466.33334
$ws.Range("M18").Value = -265.5  # was -182.33334

$ws.Range("H74").Value = 4059.4  # was 4059.6
$ws.Range("I74").Value = 3432.3333  # was 4149
$ws.Range("J74").Value = 5000  # was 4000
$ws.Range("K74").Value = 3432.3333  # was 4149
$ws.Range("L74").Value = 5000  # was 4000
$ws.Range("M74").Value = -2496.3333  # was -3213
$ws.Range("N74").Value = -6872  # was -5872

$ws.Range("H77").Value = 4059.4  # was 4059.6
$ws.Range("I77").Value = 3432.3333  # was 4149
$ws.Range("J77").Value = 5000  # was 4000
$ws.Range("K77").Value = 17161.6665  # was 20745
$ws.Range("L77").Value = 25000  # was 20000
$ws.Range("M77").Value = -12481.6665  # was -16065
$ws.Range("N77").Value = -34360  # was -29360

$ws.Range("H98").Value = 111167470  # was 142929500
$ws.Range("I98").Value = 125063350  # was 142929500
$ws.Range("J98").Value = 400  # was 0
$ws.Range("K98").Value = 125063350  # was 142929500
$ws.Range("L98").Value = 400  # was 0
$ws.Range("M98").Value = -125061852  # was -142928002
$ws.Range("N98").Value = -3396  # was None

$ws.Range("H114").Value = 136000  # was 132000
$ws.Range("J114").Value = 136000  # was 132000
$ws.Range("L114").Value = 136000  # was 132000
$ws.Range("N114").Value = -144678  # was -140678

$ws.Range("H122").Value = 111167470  # was 142929500
$ws.Range("I122").Value = 125063350  # was 142929500
$ws.Range("J122").Value = 400  # was 0
$ws.Range("K122").Value = 375190050  # was 428788500
$ws.Range("L122").Value = 1200  # was 0
$ws.Range("M122").Value = -375187600  # was -428786050
$ws.Range("N122").Value = -6100  # was None

$ws.Range("H126").Value = 90750  # was 74110
$ws.Range("J126").Value = 138000  # was 72775
$ws.Range("L126").Value = 138000  # was 72775
$ws.Range("N126").Value = -147880  # was -82655

$ws.Range("H135").Value = 1370.2174  # was 1478.2273
$ws.Range("I135").Value = 1427.15  # was 1514.95
$ws.Range("J135").Value = 990.6667  # was 1111
$ws.Range("K135").Value = 12844.35  # was 13634.55
$ws.Range("L135").Value = 8916.0003  # was 9999
$ws.Range("M135").Value = -10309.35  # was -11099.55
$ws.Range("N135").Value = -13986.0003  # was -15069

$ws.Range("H138").Value = 2888.681  # was 2884.6106
$ws.Range("J138").Value = 3294.7124  # was 3284
$ws.Range("L138").Value = 9884.137199999999  # was 9852
$ws.Range("N138").Value = -20164.1372  # was -20132

$ws.Range("H140").Value = 49000  # was 48999.5
$ws.Range("J140").Value = 49000  # was 48999.5
$ws.Range("L140").Value = 49000  # was 48999.5
$ws.Range("N140").Value = -59360  # was -59359.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1958  # was 2149.6
$ws.Range("J4").Value = 1999  # was 2998
$ws.Range("L4").Value = 1999  # was 2998
$ws.Range("N4").Value = -2231  # was -3230

$ws.Range("H34").Value = 0  # was 30000
$ws.Range("I34").Value = 0  # was 30000
$ws.Range("J34").Value = 0  # was 30000
$ws.Range("K34").Value = 0  # was 30000
$ws.Range("L34").Value = 0  # was 30000
$ws.Range("M34").ClearContents()  # was -29729
$ws.Range("N34").ClearContents()  # was -30542

$ws.Range("H45").Value = 3278.5386  # was 3385.0833
$ws.Range("I45").Value = 3035.6667  # was 3165.125
$ws.Range("K45").Value = 3035.6667  # was 3165.125
$ws.Range("M45").Value = -2658.6667  # was -2788.125

$ws.Range("H114").Value = 101696  # was 96522
$ws.Range("J114").Value = 101696  # was 96522
$ws.Range("L114").Value = 101696  # was 96522
$ws.Range("N114").Value = -110374  # was -105200

$ws.Range("H132").Value = 7000  # was 0
$ws.Range("I132").Value = 7000  # was 0
$ws.Range("K132").Value = 21000  # was 0
$ws.Range("M132").Value = -18470  # was None

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3165.9333  # was 3312.923
$ws.Range("I86").Value = 3001.9167  # was 3160.2
$ws.Range("K86").Value = 3001.9167  # was 3160.2
$ws.Range("M86").Value = -1878.9167  # was -2037.2

$ws.Range("H89").Value = 3165.9333  # was 3312.923
$ws.Range("I89").Value = 3001.9167  # was 3160.2
$ws.Range("K89").Value = 15009.5835  # was 15801
$ws.Range("M89").Value = -9393.583500000001  # was -10185

$ws.Range("H99").Value = 8780.723  # was 7550.2383
$ws.Range("I99").Value = 18476.334  # was 12650.667
$ws.Range("J99").Value = 3932.9167  # was 3724.9167
$ws.Range("K99").Value = 18476.334  # was 12650.667
$ws.Range("L99").Value = 3932.9167  # was 3724.9167
$ws.Range("M99").Value = -16978.334  # was -11152.667
$ws.Range("N99").Value = -6928.9167  # was -6720.9167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 5000  # was 0
$ws.Range("I23").Value = 5000  # was 0
$ws.Range("K23").Value = 5000  # was 0
$ws.Range("M23").Value = -4760  # was None

$ws.Range("H27").Value = 5000  # was 0
$ws.Range("I27").Value = 5000  # was 0
$ws.Range("K27").Value = 5000  # was 0
$ws.Range("M27").Value = -4808  # was None

$ws.Range("H105").Value = 2483  # was 2287.5386
$ws.Range("I105").Value = 2293.375  # was 2144.8333
$ws.Range("K105").Value = 2293.375  # was 2144.8333
$ws.Range("M105").Value = -546.375  # was -397.8332999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 556555  # was 1000000
$ws.Range("J9").Value = 556555  # was 1000000
$ws.Range("L9").Value = 1669665  # was 3000000
$ws.Range("N9").Value = -1670113  # was -3000448

$ws.Range("H54").Value = 16899  # was 17499.75
$ws.Range("I54").Value = 14750  # was 15000
$ws.Range("J54").Value = 18331.666  # was 18333
$ws.Range("K54").Value = 44250  # was 45000
$ws.Range("L54").Value = 54994.99800000001  # was 54999
$ws.Range("M54").Value = -43691  # was -44441
$ws.Range("N54").Value = -56112.99800000001  # was -56117

$ws.Range("H141").Value = 185173.36  # was 185173.64
$ws.Range("I141").Value = 433992.44  # was 433993.16
$ws.Range("K141").Value = 1301977.32  # was 1301979.48
$ws.Range("M141").Value = -1296797.32  # was -1296799.48

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3000  # was 2223
$ws.Range("I70").Value = 3000  # was 2169.2
$ws.Range("J70").Value = 0  # was 2492
$ws.Range("K70").Value = 3000  # was 2169.2
$ws.Range("L70").Value = 0  # was 2492
$ws.Range("M70").Value = -2730  # was -1899.2
$ws.Range("N70").ClearContents()  # was -3032

$ws.Range("H73").Value = 3000  # was 2223
$ws.Range("I73").Value = 3000  # was 2169.2
$ws.Range("J73").Value = 0  # was 2492
$ws.Range("K73").Value = 3000  # was 2169.2
$ws.Range("L73").Value = 0  # was 2492
$ws.Range("M73").Value = -2064  # was -1233.2
$ws.Range("N73").ClearContents()  # was -4364

$ws.Range("H80").Value = 4197.6  # was 4249.75
$ws.Range("I80").Value = 3997  # was 3999.6667
$ws.Range("K80").Value = 3997  # was 3999.6667
$ws.Range("M80").Value = -2999  # was -3001.6667

$ws.Range("H83").Value = 4197.6  # was 4249.75
$ws.Range("I83").Value = 3997  # was 3999.6667
$ws.Range("K83").Value = 19985  # was 19998.3335
$ws.Range("M83").Value = -14993  # was -15006.3335

$ws.Range("H126").Value = 4647.25  # was 4797.5
$ws.Range("I126").Value = 4598  # was 4800
$ws.Range("K126").Value = 13794  # was 14400
$ws.Range("M126").Value = -11324  # was -11930

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1478.9445  # was 1362.8
$ws.Range("I16").Value = 1288.8125  # was 1180.8889
$ws.Range("K16").Value = 1288.8125  # was 1180.8889
$ws.Range("M16").Value = -1118.8125  # was -1010.8889

$ws.Range("H22").Value = 619  # was 633.875
$ws.Range("I22").Value = 546.375  # was 553
$ws.Range("K22").Value = 546.375  # was 553
$ws.Range("M22").Value = -251.375  # was -258

$ws.Range("H27").Value = 619  # was 633.875
$ws.Range("I27").Value = 546.375  # was 553
$ws.Range("K27").Value = 546.375  # was 553
$ws.Range("M27").Value = -439.375  # was -446

$ws.Range("H46").Value = 4435.88  # was 4161.074
$ws.Range("I46").Value = 1799.8572  # was 1733.2
$ws.Range("J46").Value = 7790.8184  # was 7195.9165
$ws.Range("K46").Value = 1799.8572  # was 1733.2
$ws.Range("L46").Value = 7790.8184  # was 7195.9165
$ws.Range("M46").Value = -1611.8572  # was -1545.2
$ws.Range("N46").Value = -8166.8184  # was -7571.9165

$ws.Range("H55").Value = 76923670  # was 62500572
$ws.Range("I55").Value = 100000530  # was 83333850
$ws.Range("J55").Value = 827.3333  # was 740.25
$ws.Range("K55").Value = 100000530  # was 83333850
$ws.Range("L55").Value = 827.3333  # was 740.25
$ws.Range("M55").Value = -100000357  # was -83333677
$ws.Range("N55").Value = -1173.3333  # was -1086.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 109428.14  # was 109999.57
$ws.Range("J2").Value = 32999.5  # was 34999.5
$ws.Range("L2").Value = 32999.5  # was 34999.5
$ws.Range("N2").Value = -33223.5  # was -35223.5

$ws.Range("H4").Value = 73223.3  # was 68339.37
$ws.Range("J4").Value = 40318.715  # was 37716.375
$ws.Range("L4").Value = 40318.715  # was 37716.375
$ws.Range("N4").Value = -40544.715  # was -37942.375

$ws.Range("H5").Value = 3338146.2  # was 2862482.5
$ws.Range("J5").Value = 4005775.5  # was 3339562.8
$ws.Range("L5").Value = 4005775.5  # was 3339562.8
$ws.Range("N5").Value = -4005999.5  # was -3339786.8

$ws.Range("H96").Value = 1040.4  # was 1033.5
$ws.Range("I96").Value = 1040.4  # was 1033.5
$ws.Range("K96").Value = 1040.4  # was 1033.5
$ws.Range("M96").Value = 332.5999999999999  # was 339.5
